# Houston Rockets 2025-26 stats: add Rebounds, 3PM, Avg Rebounds, Avg 3PM sheets
$wb = $excel.ActiveWorkbook

# --- Create "Rebounds" sheet as a copy of "Points" (same Game Time/Opponent columns) ---
$srcPoints = $wb.Worksheets.Item("Points")
$afterAssists = $wb.Worksheets.Item("Assists")
$srcPoints.Copy($null, $afterAssists)
$wsReb = $wb.Worksheets.Item("Points (2)")
$wsReb.Name = "Rebounds"

# Overwrite the per-player stat columns (C:O) with rebounds numbers
$wsReb.Range("C2").Value = 0
$wsReb.Range("D2").Value = 4
$wsReb.Range("E2").Value = 0
$wsReb.Range("F2").Value = 9
$wsReb.Range("G2").Value = 0
$wsReb.Range("H2").Value = 5
$wsReb.Range("I2").Value = 13
$wsReb.Range("J2").Value = 2
$wsReb.Range("K2").Value = 6
$wsReb.Range("L2").Value = 2
$wsReb.Range("M2").Value = 11
$wsReb.Range("N2").Value = 0
$wsReb.Range("O2").Value = 0
$wsReb.Range("C3").Value = 0
$wsReb.Range("D3").Value = 4
$wsReb.Range("E3").Value = 0
$wsReb.Range("F3").Value = 3
$wsReb.Range("G3").Value = 0
$wsReb.Range("H3").Value = 6
$wsReb.Range("I3").Value = 10
$wsReb.Range("J3").Value = 2
$wsReb.Range("K3").Value = 5
$wsReb.Range("L3").Value = 4
$wsReb.Range("M3").Value = 7
$wsReb.Range("N3").Value = 4
$wsReb.Range("O3").Value = 0
$wsReb.Range("C4").Value = 0
$wsReb.Range("D4").Value = 4
$wsReb.Range("E4").Value = 0
$wsReb.Range("F4").Value = 6
$wsReb.Range("G4").Value = 1
$wsReb.Range("H4").Value = 4
$wsReb.Range("I4").Value = 8
$wsReb.Range("J4").Value = 4
$wsReb.Range("K4").Value = 5
$wsReb.Range("L4").Value = 0
$wsReb.Range("M4").Value = 6
$wsReb.Range("N4").Value = 1
$wsReb.Range("O4").Value = 0
$wsReb.Range("C5").Value = 0
$wsReb.Range("D5").Value = 9
$wsReb.Range("E5").Value = 0
$wsReb.Range("F5").Value = 5
$wsReb.Range("G5").Value = 0
$wsReb.Range("H5").Value = 5
$wsReb.Range("I5").Value = 12
$wsReb.Range("J5").Value = 5
$wsReb.Range("K5").Value = 5
$wsReb.Range("L5").Value = 1
$wsReb.Range("M5").Value = 8
$wsReb.Range("N5").Value = 3
$wsReb.Range("O5").Value = 0
$wsReb.Range("C6").Value = 2
$wsReb.Range("D6").Value = 9
$wsReb.Range("E6").Value = 0
$wsReb.Range("F6").Value = 4
$wsReb.Range("G6").Value = 1
$wsReb.Range("H6").Value = 6
$wsReb.Range("I6").Value = 5
$wsReb.Range("J6").Value = 0
$wsReb.Range("K6").Value = 6
$wsReb.Range("L6").Value = 0
$wsReb.Range("M6").Value = 10
$wsReb.Range("N6").Value = 10
$wsReb.Range("O6").Value = 0
$wsReb.Range("C7").Value = 1
$wsReb.Range("D7").Value = 5
$wsReb.Range("E7").Value = 0
$wsReb.Range("F7").Value = 6
$wsReb.Range("G7").Value = 2
$wsReb.Range("H7").Value = 0
$wsReb.Range("I7").Value = 9
$wsReb.Range("J7").Value = 5
$wsReb.Range("K7").Value = 4
$wsReb.Range("L7").Value = 3
$wsReb.Range("M7").Value = 11
$wsReb.Range("N7").Value = 8
$wsReb.Range("O7").Value = 0
$wsReb.Range("C8").Value = 0
$wsReb.Range("D8").Value = 10
$wsReb.Range("E8").Value = 0
$wsReb.Range("F8").Value = 4
$wsReb.Range("G8").Value = 0
$wsReb.Range("H8").Value = 4
$wsReb.Range("I8").Value = 0
$wsReb.Range("J8").Value = 1
$wsReb.Range("K8").Value = 5
$wsReb.Range("L8").Value = 3
$wsReb.Range("M8").Value = 16
$wsReb.Range("N8").Value = 10
$wsReb.Range("O8").Value = 1
$wsReb.Range("C9").Value = 0
$wsReb.Range("D9").Value = 7
$wsReb.Range("E9").Value = 0
$wsReb.Range("F9").Value = 1
$wsReb.Range("G9").Value = 0
$wsReb.Range("H9").Value = 8
$wsReb.Range("I9").Value = 8
$wsReb.Range("J9").Value = 0
$wsReb.Range("K9").Value = 2
$wsReb.Range("L9").Value = 3
$wsReb.Range("M9").Value = 9
$wsReb.Range("N9").Value = 3
$wsReb.Range("O9").Value = 0
$wsReb.Range("C10").Value = 0
$wsReb.Range("D10").Value = 8
$wsReb.Range("E10").Value = 0
$wsReb.Range("F10").Value = 3
$wsReb.Range("G10").Value = 0
$wsReb.Range("H10").Value = 9
$wsReb.Range("I10").Value = 6
$wsReb.Range("J10").Value = 0
$wsReb.Range("K10").Value = 8
$wsReb.Range("L10").Value = 1
$wsReb.Range("M10").Value = 11
$wsReb.Range("N10").Value = 4
$wsReb.Range("O10").Value = 0

# --- Create "3PM" sheet as a copy of "Points", placed after "Rebounds" ---
$afterRebounds = $wb.Worksheets.Item("Rebounds")
$srcPoints.Copy($null, $afterRebounds)
$ws3pm = $wb.Worksheets.Item("Points (2)")
$ws3pm.Name = "3PM"

# Overwrite the per-player stat columns (C:O) with 3-pointers made numbers
$ws3pm.Range("C2").Value = 0
$ws3pm.Range("D2").Value = 0
$ws3pm.Range("E2").Value = 0
$ws3pm.Range("F2").Value = 0
$ws3pm.Range("G2").Value = 0
$ws3pm.Range("H2").Value = 2
$ws3pm.Range("I2").Value = 0
$ws3pm.Range("J2").Value = 2
$ws3pm.Range("K2").Value = 1
$ws3pm.Range("L2").Value = 1
$ws3pm.Range("M2").Value = 5
$ws3pm.Range("N2").Value = 0
$ws3pm.Range("O2").Value = 0
$ws3pm.Range("C3").Value = 0
$ws3pm.Range("D3").Value = 0
$ws3pm.Range("E3").Value = 0
$ws3pm.Range("F3").Value = 3
$ws3pm.Range("G3").Value = 0
$ws3pm.Range("H3").Value = 3
$ws3pm.Range("I3").Value = 0
$ws3pm.Range("J3").Value = 3
$ws3pm.Range("K3").Value = 0
$ws3pm.Range("L3").Value = 1
$ws3pm.Range("M3").Value = 0
$ws3pm.Range("N3").Value = 0
$ws3pm.Range("O3").Value = 0
$ws3pm.Range("C4").Value = 1
$ws3pm.Range("D4").Value = 0
$ws3pm.Range("E4").Value = 1
$ws3pm.Range("F4").Value = 0
$ws3pm.Range("G4").Value = 0
$ws3pm.Range("H4").Value = 1
$ws3pm.Range("I4").Value = 0
$ws3pm.Range("J4").Value = 3
$ws3pm.Range("K4").Value = 5
$ws3pm.Range("L4").Value = 2
$ws3pm.Range("M4").Value = 2
$ws3pm.Range("N4").Value = 0
$ws3pm.Range("O4").Value = 1
$ws3pm.Range("C5").Value = 0
$ws3pm.Range("D5").Value = 0
$ws3pm.Range("E5").Value = 0
$ws3pm.Range("F5").Value = 4
$ws3pm.Range("G5").Value = 0
$ws3pm.Range("H5").Value = 4
$ws3pm.Range("I5").Value = 0
$ws3pm.Range("J5").Value = 0
$ws3pm.Range("K5").Value = 2
$ws3pm.Range("L5").Value = 2
$ws3pm.Range("M5").Value = 1
$ws3pm.Range("N5").Value = 0
$ws3pm.Range("O5").Value = 0
$ws3pm.Range("C6").Value = 0
$ws3pm.Range("D6").Value = 1
$ws3pm.Range("E6").Value = 1
$ws3pm.Range("F6").Value = 2
$ws3pm.Range("G6").Value = 0
$ws3pm.Range("H6").Value = 1
$ws3pm.Range("I6").Value = 0
$ws3pm.Range("J6").Value = 4
$ws3pm.Range("K6").Value = 4
$ws3pm.Range("L6").Value = 3
$ws3pm.Range("M6").Value = 1
$ws3pm.Range("N6").Value = 0
$ws3pm.Range("O6").Value = 2
$ws3pm.Range("C7").Value = 1
$ws3pm.Range("D7").Value = 1
$ws3pm.Range("E7").Value = 0
$ws3pm.Range("F7").Value = 2
$ws3pm.Range("G7").Value = 0
$ws3pm.Range("H7").Value = 0
$ws3pm.Range("I7").Value = 0
$ws3pm.Range("J7").Value = 1
$ws3pm.Range("K7").Value = 3
$ws3pm.Range("L7").Value = 0
$ws3pm.Range("M7").Value = 0
$ws3pm.Range("N7").Value = 0
$ws3pm.Range("O7").Value = 0
$ws3pm.Range("C8").Value = 0
$ws3pm.Range("D8").Value = 1
$ws3pm.Range("E8").Value = 0
$ws3pm.Range("F8").Value = 1
$ws3pm.Range("G8").Value = 0
$ws3pm.Range("H8").Value = 2
$ws3pm.Range("I8").Value = 0
$ws3pm.Range("J8").Value = 2
$ws3pm.Range("K8").Value = 4
$ws3pm.Range("L8").Value = 3
$ws3pm.Range("M8").Value = 0
$ws3pm.Range("N8").Value = 0
$ws3pm.Range("O8").Value = 0
$ws3pm.Range("C9").Value = 0
$ws3pm.Range("D9").Value = 2
$ws3pm.Range("E9").Value = 0
$ws3pm.Range("F9").Value = 1
$ws3pm.Range("G9").Value = 0
$ws3pm.Range("H9").Value = 0
$ws3pm.Range("I9").Value = 0
$ws3pm.Range("J9").Value = 4
$ws3pm.Range("K9").Value = 3
$ws3pm.Range("L9").Value = 2
$ws3pm.Range("M9").Value = 1
$ws3pm.Range("N9").Value = 0
$ws3pm.Range("O9").Value = 0
$ws3pm.Range("C10").Value = 0
$ws3pm.Range("D10").Value = 0
$ws3pm.Range("E10").Value = 0
$ws3pm.Range("F10").Value = 2
$ws3pm.Range("G10").Value = 0
$ws3pm.Range("H10").Value = 3
$ws3pm.Range("I10").Value = 0
$ws3pm.Range("J10").Value = 4
$ws3pm.Range("K10").Value = 0
$ws3pm.Range("L10").Value = 0
$ws3pm.Range("M10").Value = 1
$ws3pm.Range("N10").Value = 0
$ws3pm.Range("O10").Value = 0

# --- Create "Avg Rebounds" sheet as a copy of "Avg Points" (Player / value layout) ---
$srcAvg = $wb.Worksheets.Item("Avg Points")
$afterAvgAssists = $wb.Worksheets.Item("Avg Assists")
$srcAvg.Copy($null, $afterAvgAssists)
$wsAR = $wb.Worksheets.Item("Avg Points (2)")
$wsAR.Name = "Avg Rebounds"
$wsAR.Range("B1").Value = "Avg Rebounds"

$wsAR.Range("A2").Value = "Alperen Sengun"
$wsAR.Range("B2").Value = 9.88888888888889
$wsAR.Range("A3").Value = "Steven Adams"
$wsAR.Range("B3").Value = 8.875
$wsAR.Range("A4").Value = "Amen Thompson"
$wsAR.Range("B4").Value = 6.666666666666667
$wsAR.Range("A5").Value = "Jabari Smith Jr."
$wsAR.Range("B5").Value = 5.875
$wsAR.Range("A6").Value = "Tari Eason"
$wsAR.Range("B6").Value = 5.111111111111111
$wsAR.Range("A7").Value = "Clint Capela"
$wsAR.Range("B7").Value = 4.777777777777778
$wsAR.Range("A8").Value = "Kevin Durant"
$wsAR.Range("B8").Value = 4.555555555555555
$wsAR.Range("A9").Value = "Reed Sheppard"
$wsAR.Range("B9").Value = 2.111111111111111
$wsAR.Range("A10").Value = "Josh Okogie"
$wsAR.Range("B10").Value = 1.888888888888889
$wsAR.Range("A11").Value = "Jae'Sean Tate"
$wsAR.Range("B11").Value = 0.8
$wsAR.Range("A12").Value = "Aaron Holiday"
$wsAR.Range("B12").Value = 0.75
$wsAR.Range("A13").Value = "Jeff Green"
$wsAR.Range("B13").Value = 0.3333333333333333
$wsAR.Range("A14").Value = "JD Davison"
$wsAR.Range("B14").Value = 0

# --- Create "Avg 3PM" sheet as a copy of "Avg Points", placed after "Avg Rebounds" ---
$afterAvgRebounds = $wb.Worksheets.Item("Avg Rebounds")
$srcAvg.Copy($null, $afterAvgRebounds)
$wsA3 = $wb.Worksheets.Item("Avg Points (2)")
$wsA3.Name = "Avg 3PM"
$wsA3.Range("B1").Value = "Avg 3PM"

$wsA3.Range("A2").Value = "Reed Sheppard"
$wsA3.Range("B2").Value = 2.555555555555555
$wsA3.Range("A3").Value = "Tari Eason"
$wsA3.Range("B3").Value = 2.444444444444445
$wsA3.Range("A4").Value = "Jabari Smith Jr."
$wsA3.Range("B4").Value = 2
$wsA3.Range("A5").Value = "Kevin Durant"
$wsA3.Range("B5").Value = 1.666666666666667
$wsA3.Range("A6").Value = "Josh Okogie"
$wsA3.Range("B6").Value = 1.555555555555556
$wsA3.Range("A7").Value = "Alperen Sengun"
$wsA3.Range("B7").Value = 1.222222222222222
$wsA3.Range("A8").Value = "Jeff Green"
$wsA3.Range("B8").Value = 1
$wsA3.Range("A9").Value = "JD Davison"
$wsA3.Range("B9").Value = 0.6666666666666666
$wsA3.Range("A10").Value = "Amen Thompson"
$wsA3.Range("B10").Value = 0.5555555555555556
$wsA3.Range("A11").Value = "Aaron Holiday"
$wsA3.Range("B11").Value = 0.5
$wsA3.Range("A12").Value = "Jae'Sean Tate"
$wsA3.Range("B12").Value = 0
$wsA3.Range("A13").Value = "Steven Adams"
$wsA3.Range("B13").Value = 0
$wsA3.Range("A14").Value = "Clint Capela"
$wsA3.Range("B14").Value = 0
